# Control-signal table update:
#   - insert a new "hilowrite" column (J) between "regwrite" (I) and "jump" (old J, now K)
#   - populate the new column for every instruction row
#   - MULT / MULTU rows gain full control-signal data (previously only label cells)
#
# Columns after the insert:
#   A label, B " ", C memtoreg, D memen, E memwrite, F branch, G alusrc,
#   H regdst, I regwrite, J hilowrite (NEW), K jump, L jal, M jr, N bal, O " "

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a whole new column at J; this shifts the existing J:N data (jump/jal/jr/bal
#    and the trailing " " marker column) one slot to the right, to K:O, automatically.
$ws.Columns("J:J").Insert()

# 2. New column header.
$ws.Range("J1").Value = "hilowrite"

# 3. hilowrite = 1 only for the "move to HI/LO" and multiply instructions (MTHI, MTLO,
#    MULT, MULTU); 0 for every other populated data row; left blank for the separator
#    rows and for the immediate-type instructions (rows 32-39).
$hilowriteRows = @{
    2  = 0   # NOP
    4  = 0   # AND
    5  = 0   # OR
    6  = 0   # XOR
    7  = 0   # NOR
    9  = 0   # SLL
    10 = 0   # SRL
    11 = 0   # SRA
    12 = 0   # SLLV
    13 = 0   # SRLV
    14 = 0   # SRAV
    16 = 0   # MFHI
    17 = 0   # MFLO
    18 = 1   # MTHI
    19 = 1   # MTLO
    21 = 0   # ADD
    22 = 0   # ADDU
    23 = 0   # SUB
    24 = 0   # SUBU
    25 = 0   # SLT
    26 = 0   # SLTU
    27 = 1   # MULT
    28 = 1   # MULTU
}

foreach ($r in $hilowriteRows.Keys) {
    $ws.Cells.Item($r, 10).Value = $hilowriteRows[$r]
}

# 4. MULT (row 27) and MULTU (row 28) previously only had the label (A) and the
#    blank-marker (B) cells filled in; now they get the full control-signal row,
#    matching the other R-type ALU rows (regdst=1, everything else but hilowrite 0).
foreach ($r in 27, 28) {
    $ws.Cells.Item($r, 3).Value  = 0   # C memtoreg
    $ws.Cells.Item($r, 4).Value  = 0   # D memen
    $ws.Cells.Item($r, 5).Value  = 0   # E memwrite
    $ws.Cells.Item($r, 6).Value  = 0   # F branch
    $ws.Cells.Item($r, 7).Value  = 0   # G alusrc
    $ws.Cells.Item($r, 8).Value  = 1   # H regdst
    $ws.Cells.Item($r, 9).Value  = 0   # I regwrite
    # J (hilowrite) already set to 1 above
    $ws.Cells.Item($r, 11).Value = 0   # K jump
    $ws.Cells.Item($r, 12).Value = 0   # L jal
    $ws.Cells.Item($r, 13).Value = 0   # M jr
    $ws.Cells.Item($r, 14).Value = 0   # N bal
}
